$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels (shared strings) to simpler, underscore-joined names
$ws.Range("B1").Value = "Total"
$ws.Range("C1").Value = "In_county_of_residence"
$ws.Range("D1").Value = "Outside_county_of_residence"

# Move the active selection from C9 to D9
[void]$ws.Range("D9").Select()
